$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 16 with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(16, 1).Value = "Mijn retour is nog steeds niet verwerkt."
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$logs.Cells.Item(16, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(16, 5).Value = "Beste klant,`nBedankt voor uw bericht. Ik begrijp uw zorgen over de verwerking van uw retourzending. Om u beter van dienst te kunnen zijn, heb ik wat meer informatie nodig. Kunt u alstublieft het volgende verstrekken:`n1. Het ordernummer van uw aankoop.`n2. De datum waarop u de retour heeft verzonden.`n3. Eventuele traceerinformatie van de retourzending.`nMet deze gegevens kunnen we verder onderzoek doen naar de status van uw retour en u zo snel mogelijk een update geven.`nIk kijk uit naar uw reactie.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Cells.Item(16, 6).Value = "2025-08-01 23:58:23"
$logs.Cells.Item(16, 7).Value = "Ja"
$logs.Cells.Item(16, 8).Value = "Nee"
$logs.Cells.Item(16, 9).Value = "Ja"
$logs.Cells.Item(16, 10).Value = "Nee"

# The multi-line text in column E triggers an automatic row-height change;
# AutoFit restores the row to use the sheet's default (non-custom) height.
$logs.Rows.Item(16).EntireRow.AutoFit()

# Extend the conditional-formatting ranges so row 16 is covered too
$logs.Range("D2:D15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D16"))
$logs.Range("G2:G15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G16"))
$logs.Range("H2:H15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H16"))
$logs.Range("I2:I15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I16"))
$logs.Range("J2:J15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J16"))

# --- Sheet "Dashboard": add the new category count row ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(6, 2).Value = 1

# --- Chart on Dashboard: extend the series ranges to include row 6 ---
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
